# Update the "F" column (view/favorite counters) across all four sheets
# to reflect newly scraped numbers, per commit "Update gh-pages to output
# generated at 456a3b4".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 331
$ws.Range("F5").Value = 2187
$ws.Range("F10").Value = 1459
$ws.Range("F11").Value = 1773
$ws.Range("F14").Value = 625
$ws.Range("F16").Value = 2252
$ws.Range("F17").Value = 1263
$ws.Range("F20").Value = 1873
$ws.Range("F22").Value = 5590
$ws.Range("F23").Value = 1058
$ws.Range("F25").Value = 90
$ws.Range("F26").Value = 1345
$ws.Range("F27").Value = 246
$ws.Range("F29").Value = 556
$ws.Range("F30").Value = 127
$ws.Range("F33").Value = 1182
$ws.Range("F35").Value = 3524
$ws.Range("F37").Value = 1131
$ws.Range("F44").Value = 813
$ws.Range("F46").Value = 39
$ws.Range("F47").Value = 4
$ws.Range("F49").Value = 33
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 148658
$ws.Range("F18").Value = 123
$ws.Range("F36").Value = 13
$ws.Range("F43").Value = 149
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 1412
$ws.Range("F11").Value = 2636
$ws.Range("F12").Value = 190
$ws.Range("F13").Value = 322
$ws.Range("F14").Value = 1018
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 331
$ws.Range("F8").Value = 2636
$ws.Range("F9").Value = 2187
$ws.Range("F14").Value = 1459
$ws.Range("F16").Value = 1773
$ws.Range("F18").Value = 625
$ws.Range("F21").Value = 2252
$ws.Range("F22").Value = 190
$ws.Range("F23").Value = 1263
$ws.Range("F26").Value = 5590
$ws.Range("F27").Value = 322
$ws.Range("F28").Value = 1058
$ws.Range("F30").Value = 1018
$ws.Range("F31").Value = 1345
$ws.Range("F35").Value = 556
$ws.Range("F36").Value = 127
$ws.Range("F38").Value = 1182
$ws.Range("F39").Value = 3524
$ws.Range("F42").Value = 1131
$ws.Range("F47").Value = 813
$ws.Range("F49").Value = 149
$ws.Range("F50").Value = 149
